# Applies the MPA test automation upload update to the "Data" worksheet.
# Rows 7, 11, 12, 16, 17, 21, 22, 26, 27 (columns A:AB) are replaced with
# the single numeric value 60000169.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$rows = @(7, 11, 12, 16, 17, 21, 22, 26, 27)

foreach ($r in $rows) {
    $ws.Range("A$r`:AB$r").Value = 60000169
}
